$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.178.81'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.615.59'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.570'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.35%  '
$ws.Range("D9").Value = '2.618.56'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '3.074.42'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '60.208.44'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.50%  '
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("D18").Value = '2.618.30'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("E28").Value = '  -4.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.88%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("E32").Value = '  -0.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.921'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.862'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  -3.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '289.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.623'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0547'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("E50").Value = '  -2.45%  '
$ws.Range("D51").Value = '1.954.20'
$ws.Range("E51").Value = '  -0.95%  '
